$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 278, shifting existing rows 278:292 down to 279:293.
$ws.Rows.Item(278).Insert()

# Populate the newly inserted row 278 with the new weekly data point.
$ws.Cells.Item(278, 1).Value = 3
$ws.Cells.Item(278, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(278, 3).Value = "Coquimbo"
$ws.Cells.Item(278, 4).Value = 45267
$ws.Cells.Item(278, 5).Value = 5
$ws.Cells.Item(278, 6).Value = 100112030
$ws.Cells.Item(278, 7).Value = "Poroto granado"
$ws.Cells.Item(278, 8).Value = "Sin especificar"
$ws.Cells.Item(278, 9).Value = "Primera"
$ws.Cells.Item(278, 10).Value = 30
$ws.Cells.Item(278, 11).Value = 33000
$ws.Cells.Item(278, 12).Value = 33000
$ws.Cells.Item(278, 13).Value = 33000
$ws.Cells.Item(278, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(278, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(278, 16).Value = 1320
$ws.Cells.Item(278, 17).Value = 25
$ws.Cells.Item(278, 18).Value = "Hortaliza"
